$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the new, longer description text
# (the COM width<->stored-units conversion adds a sub-character offset, so
# 27.15 is the input that lands exactly on a stored width of 28)
$ws.Columns("B").ColumnWidth = 27.15

# New BOM row: 5M Silver Plated Copper Wire 24AWG
$ws.Range("F6").Value = "https://www.aliexpress.com/item/1005005321856209.html"
$ws.Range("A6").Value = "5M Silver Plated Copper Wire 24AWG"
$ws.Range("B6").Value = "for connecting esp32 and led strips"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 3.46
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("G6").Formula = "=ROUND(E6*1.12,2)+G5"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.aliexpress.com/item/1005005321856209.html") | Out-Null

# Hyperlinks.Add applies its own font formatting; reapply the same
# "Hyperlink" cell style already used by the other URL cells (F2:F5) so F6
# matches them exactly instead of minting a new style variant
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to D7, matching the post-edit saved state
$ws.Range("D7").Select()
